$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 18,20
$arr[0,0] = 'ECs'
$arr[0,1] = 'Cirbp'
$arr[0,2] = 'Trem1'
$arr[0,3] = 'Inflammatory-Mac'
$arr[0,4] = 2
$arr[0,5] = 1
$arr[0,6] = 17.15177
$arr[0,7] = 34.30354
$arr[0,8] = 0.3417885598906996
$arr[0,9] = 0.2838492023609646
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 12.97650433333333
$arr[0,13] = 38.929513
$arr[0,14] = 0.06907304296508091
$arr[0,15] = 0.06907304296508089
$arr[0,16] = 222.5700177293367
$arr[0,17] = 1335.42010637602
$arr[0,18] = 0.02360837588230342
$arr[0,19] = 0.01960632815028284
$arr[1,0] = 'ECs'
$arr[1,1] = 'Cirbp'
$arr[1,2] = 'Trem1'
$arr[1,3] = 'Neutrophils'
$arr[1,4] = 2
$arr[1,5] = 1
$arr[1,6] = 17.15177
$arr[1,7] = 34.30354
$arr[1,8] = 0.3417885598906996
$arr[1,9] = 0.2838492023609646
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 172.5923843333333
$arr[1,13] = 517.777153
$arr[1,14] = 0.918697429775355
$arr[1,15] = 0.918697429775355
$arr[1,16] = 2960.264879836937
$arr[1,17] = 17761.58927902162
$arr[1,18] = 0.3140002714982056
$arr[1,19] = 0.2607715326528027
$arr[2,0] = 'ECs'
$arr[2,1] = 'Cirbp'
$arr[2,2] = 'Trem1'
$arr[2,3] = 'Resolving-Mac'
$arr[2,4] = 2
$arr[2,5] = 1
$arr[2,6] = 17.15177
$arr[2,7] = 34.30354
$arr[2,8] = 0.3417885598906996
$arr[2,9] = 0.2838492023609646
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 2.297517333333333
$arr[2,13] = 6.892552
$arr[2,14] = 0.01222952725956408
$arr[2,15] = 0.01222952725956408
$arr[2,16] = 39.40648887234666
$arr[2,17] = 236.43893323408
$arr[2,18] = 0.004179912510190461
$arr[2,19] = 0.003471341557878938
$arr[3,0] = 'FAPs'
$arr[3,1] = 'Cirbp'
$arr[3,2] = 'Trem1'
$arr[3,3] = 'Inflammatory-Mac'
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 4.783142666666667
$arr[3,7] = 14.349428
$arr[3,8] = 0.09531514495540462
$arr[3,9] = 0.1187362497321294
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 12.97650433333333
$arr[3,13] = 38.929513
$arr[3,14] = 0.06907304296508091
$arr[3,15] = 0.06907304296508089
$arr[3,16] = 62.06847154095155
$arr[3,17] = 558.616243868564
$arr[3,18] = 0.006583707102727578
$arr[3,19] = 0.008201474079259951
$arr[4,0] = 'FAPs'
$arr[4,1] = 'Cirbp'
$arr[4,2] = 'Trem1'
$arr[4,3] = 'Neutrophils'
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 4.783142666666667
$arr[4,7] = 14.349428
$arr[4,8] = 0.09531514495540462
$arr[4,9] = 0.1187362497321294
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 172.5923843333333
$arr[4,13] = 517.777153
$arr[4,14] = 0.918697429775355
$arr[4,15] = 0.918697429775355
$arr[4,16] = 825.5339974464982
$arr[4,17] = 7429.805977018484
$arr[4,18] = 0.08756577868919561
$arr[4,19] = 0.109082687450072
$arr[5,0] = 'FAPs'
$arr[5,1] = 'Cirbp'
$arr[5,2] = 'Trem1'
$arr[5,3] = 'Resolving-Mac'
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 4.783142666666667
$arr[5,7] = 14.349428
$arr[5,8] = 0.09531514495540462
$arr[5,9] = 0.1187362497321294
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 2.297517333333333
$arr[5,13] = 6.892552
$arr[5,14] = 0.01222952725956408
$arr[5,15] = 0.01222952725956408
$arr[5,16] = 10.98935318447289
$arr[5,17] = 98.90417866025599
$arr[5,18] = 0.001165659163481423
$arr[5,19] = 0.001452088202797485
$arr[6,0] = 'Inflammatory-Mac'
$arr[6,1] = 'Cirbp'
$arr[6,2] = 'Trem1'
$arr[6,3] = 'Inflammatory-Mac'
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 3.755572666666666
$arr[6,7] = 11.266718
$arr[6,8] = 0.07483844368860322
$arr[6,9] = 0.09322795599305267
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 12.97650433333333
$arr[6,13] = 38.929513
$arr[6,14] = 0.06907304296508091
$arr[6,15] = 0.06907304296508089
$arr[6,16] = 48.73420498314822
$arr[6,17] = 438.607844848334
$arr[6,18] = 0.005169319036342678
$arr[6,19] = 0.006439538609854798
$arr[7,0] = 'Inflammatory-Mac'
$arr[7,1] = 'Cirbp'
$arr[7,2] = 'Trem1'
$arr[7,3] = 'Neutrophils'
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 3.755572666666666
$arr[7,7] = 11.266718
$arr[7,8] = 0.07483844368860322
$arr[7,9] = 0.09322795599305267
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 172.5923843333333
$arr[7,13] = 517.777153
$arr[7,14] = 0.918697429775355
$arr[7,15] = 0.918697429775355
$arr[7,16] = 648.1832410770949
$arr[7,17] = 5833.649169693854
$arr[7,18] = 0.06875388586510742
$arr[7,19] = 0.08564828355402739
$arr[8,0] = 'Inflammatory-Mac'
$arr[8,1] = 'Cirbp'
$arr[8,2] = 'Trem1'
$arr[8,3] = 'Resolving-Mac'
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 3.755572666666666
$arr[8,7] = 11.266718
$arr[8,8] = 0.07483844368860322
$arr[8,9] = 0.09322795599305267
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 2.297517333333333
$arr[8,13] = 6.892552
$arr[8,14] = 0.01222952725956408
$arr[8,15] = 0.01222952725956408
$arr[8,16] = 8.628493298259555
$arr[8,17] = 77.65643968433599
$arr[8,18] = 0.0009152387871531247
$arr[8,19] = 0.001140133829170478
$arr[9,0] = 'MuSCs'
$arr[9,1] = 'Cirbp'
$arr[9,2] = 'Trem1'
$arr[9,3] = 'Inflammatory-Mac'
$arr[9,4] = 2
$arr[9,5] = 1
$arr[9,6] = 12.544146
$arr[9,7] = 25.088292
$arr[9,8] = 0.2499710290190855
$arr[9,9] = 0.2075964076243726
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 12.97650433333333
$arr[9,13] = 38.929513
$arr[9,14] = 0.06907304296508091
$arr[9,15] = 0.06907304296508089
$arr[9,16] = 162.779164926966
$arr[9,17] = 976.6749895617961
$arr[9,18] = 0.01726625962746078
$arr[9,19] = 0.01433931558323473
$arr[10,0] = 'MuSCs'
$arr[10,1] = 'Cirbp'
$arr[10,2] = 'Trem1'
$arr[10,3] = 'Neutrophils'
$arr[10,4] = 2
$arr[10,5] = 1
$arr[10,6] = 12.544146
$arr[10,7] = 25.088292
$arr[10,8] = 0.2499710290190855
$arr[10,9] = 0.2075964076243726
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 172.5923843333333
$arr[10,13] = 517.777153
$arr[10,14] = 0.918697429775355
$arr[10,15] = 0.918697429775355
$arr[10,16] = 2165.024067565446
$arr[10,17] = 12990.14440539268
$arr[10,18] = 0.2296477418781345
$arr[10,19] = 0.190718286115108
$arr[11,0] = 'MuSCs'
$arr[11,1] = 'Cirbp'
$arr[11,2] = 'Trem1'
$arr[11,3] = 'Resolving-Mac'
$arr[11,4] = 2
$arr[11,5] = 1
$arr[11,6] = 12.544146
$arr[11,7] = 25.088292
$arr[11,8] = 0.2499710290190855
$arr[11,9] = 0.2075964076243726
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 2.297517333333333
$arr[11,13] = 6.892552
$arr[11,14] = 0.01222952725956408
$arr[11,15] = 0.01222952725956408
$arr[11,16] = 28.820392866864
$arr[11,17] = 172.922357201184
$arr[11,18] = 0.00305702751349019
$arr[11,19] = 0.002538805926029841
$arr[12,0] = 'Neutrophils'
$arr[12,1] = 'Cirbp'
$arr[12,2] = 'Trem1'
$arr[12,3] = 'Inflammatory-Mac'
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 8.593081
$arr[12,7] = 25.779243
$arr[12,8] = 0.1712369498899608
$arr[12,9] = 0.2133137735353109
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 12.97650433333333
$arr[12,13] = 38.929513
$arr[12,14] = 0.06907304296508091
$arr[12,15] = 0.06907304296508089
$arr[12,16] = 111.5081528331843
$arr[12,17] = 1003.573375498659
$arr[12,18] = 0.01182785719695867
$arr[12,19] = 0.01473423144444807
$arr[13,0] = 'Neutrophils'
$arr[13,1] = 'Cirbp'
$arr[13,2] = 'Trem1'
$arr[13,3] = 'Neutrophils'
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 8.593081
$arr[13,7] = 25.779243
$arr[13,8] = 0.1712369498899608
$arr[13,9] = 0.2133137735353109
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 172.5923843333333
$arr[13,13] = 517.777153
$arr[13,14] = 0.918697429775355
$arr[13,15] = 0.918697429775355
$arr[13,16] = 1483.100338559464
$arr[13,17] = 13347.90304703518
$arr[13,18] = 0.1573149457464782
$arr[13,19] = 0.1959708154825723
$arr[14,0] = 'Neutrophils'
$arr[14,1] = 'Cirbp'
$arr[14,2] = 'Trem1'
$arr[14,3] = 'Resolving-Mac'
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 8.593081
$arr[14,7] = 25.779243
$arr[14,8] = 0.1712369498899608
$arr[14,9] = 0.2133137735353109
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 2.297517333333333
$arr[14,13] = 6.892552
$arr[14,14] = 0.01222952725956408
$arr[14,15] = 0.01222952725956408
$arr[14,16] = 19.74275254423733
$arr[14,17] = 177.684772898136
$arr[14,18] = 0.002094146946523884
$arr[14,19] = 0.002608726608290565
$arr[15,0] = 'Resolving-Mac'
$arr[15,1] = 'Cirbp'
$arr[15,2] = 'Trem1'
$arr[15,3] = 'Inflammatory-Mac'
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 3.354687
$arr[15,7] = 10.064061
$arr[15,8] = 0.06684987255624644
$arr[15,9] = 0.08327641075416972
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 12.97650433333333
$arr[15,13] = 38.929513
$arr[15,14] = 0.06907304296508091
$arr[15,15] = 0.06907304296508089
$arr[15,16] = 43.53211039247699
$arr[15,17] = 391.788993532293
$arr[15,18] = 0.004617524119287794
$arr[15,19] = 0.00575215509800049
$arr[16,0] = 'Resolving-Mac'
$arr[16,1] = 'Cirbp'
$arr[16,2] = 'Trem1'
$arr[16,3] = 'Neutrophils'
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 3.354687
$arr[16,7] = 10.064061
$arr[16,8] = 0.06684987255624644
$arr[16,9] = 0.08327641075416972
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 172.5923843333333
$arr[16,13] = 517.777153
$arr[16,14] = 0.918697429775355
$arr[16,15] = 0.918697429775355
$arr[16,16] = 578.9934280220369
$arr[16,17] = 5210.940852198332
$arr[16,18] = 0.06141480609823365
$arr[16,19] = 0.07650582452077245
$arr[17,0] = 'Resolving-Mac'
$arr[17,1] = 'Cirbp'
$arr[17,2] = 'Trem1'
$arr[17,3] = 'Resolving-Mac'
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 3.354687
$arr[17,7] = 10.064061
$arr[17,8] = 0.06684987255624644
$arr[17,9] = 0.08327641075416972
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 2.297517333333333
$arr[17,13] = 6.892552
$arr[17,14] = 0.01222952725956408
$arr[17,15] = 0.01222952725956408
$arr[17,16] = 7.707451530408
$arr[17,17] = 69.367063773672
$arr[17,18] = 0.0008175423387250007
$arr[17,19] = 0.001018431135396774

$ws.Range("A2:T19").Value = $arr
